$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use D25 (untouched, General-format style) as a stable formatting template
# so assigning date-look-alike strings into column E keeps the original
# numeric cell style instead of Excel auto-converting them to date serials.
$formatSource = $ws.Range("D25")

# Row 13
$ws.Range("A13").Value = "Intel(R) Wi-Fi 6 AX200 160MHz - 23.70.2.3"
$ws.Range("B13").Value = 18721
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "2024-07-23"
$formatSource.Copy()
$ws.Range("E13").PasteSpecial(-4122)

# Row 14
$ws.Range("A14").Value = "Intel(R) Wi-Fi 6 AX200 160MHz - 22.250.10.1"
$ws.Range("B14").Value = 69578
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "2023-08-14"
$formatSource.Copy()
$ws.Range("E14").PasteSpecial(-4122)

# Row 15
$ws.Range("A15").Value = "Intel(R) Wi-Fi 6 AX200 160MHz - 22.230.0.8"
$ws.Range("B15").Value = 329845
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "2023-05-08"
$formatSource.Copy()
$ws.Range("E15").PasteSpecial(-4122)

# Row 16
$ws.Range("A16").Value = "Intel(R) Wi-Fi 6 AX200 160MHz - 22.200.0.6"
$ws.Range("B16").Value = 143808
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "2023-01-16"
$formatSource.Copy()
$ws.Range("E16").PasteSpecial(-4122)

# Row 17
$ws.Range("A17").Value = "Intel(R) Wi-Fi 6 AX200 160MHz - 22.190.0.4"
$ws.Range("B17").Value = 287148
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "2022-11-22"
$formatSource.Copy()
$ws.Range("E17").PasteSpecial(-4122)

# Row 18
$ws.Range("A18").Value = "Intel(R) Wi-Fi 6 AX200 160MHz - 22.160.0.4"
$ws.Range("B18").Value = 96526
$ws.Range("D18").Value = 99.9
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "2022-08-13"
$formatSource.Copy()
$ws.Range("E18").PasteSpecial(-4122)

# Row 19
$ws.Range("A19").Value = "Intel(R) Wi-Fi 6 AX200 160MHz - 22.30.0.11"
$ws.Range("B19").Value = 67111
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "2021-01-19"
$formatSource.Copy()
$ws.Range("E19").PasteSpecial(-4122)

# Row 20
$ws.Range("A20").Value = "Intel(R) Wi-Fi 6 AX200 160MHz - 22.10.0.7"
$ws.Range("B20").Value = 66577
$ws.Range("D20").Value = 100
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "2020-10-19"
$formatSource.Copy()
$ws.Range("E20").PasteSpecial(-4122)

# Row 21
$ws.Range("A21").Value = "Intel(R) Wi-Fi 6 AX200 160MHz - 22.0.1.1"
$ws.Range("B21").Value = 15734
$ws.Range("D21").Value = 99.9
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "2020-09-28"
$formatSource.Copy()
$ws.Range("E21").PasteSpecial(-4122)

# Row 22
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "2019-12-14"
$formatSource.Copy()
$ws.Range("E22").PasteSpecial(-4122)

# Row 23
$ws.Range("A23").Value = "Intel(R) Wi-Fi 6 AX200 160MHz - 21.40.2.2"
$ws.Range("B23").Value = 88435
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "2019-08-31"
$formatSource.Copy()
$ws.Range("E23").PasteSpecial(-4122)

# Row 24
$ws.Range("A24").Value = "Intel(R) Wi-Fi 6 AX200 160MHz - 21.30.4.1"
$ws.Range("B24").Value = 13016
$ws.Range("D24").Value = 100
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "2019-07-29"
$formatSource.Copy()
$ws.Range("E24").PasteSpecial(-4122)

$excel.CutCopyMode = 0